$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text values (e.g. "1.001", "24.418.50")
# that must stay as literal text, not be auto-converted to numbers by Excel.
# We temporarily force the cell format to Text ("@") while assigning the
# value, then restore the original "General" number format.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.403.66"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.666.50"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.12"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3957"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3907"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.59"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +8.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.399"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08586"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.39"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.301"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.946"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +6.22%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001343"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +4.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.663.76"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.29"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07007"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.56"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.994"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.402.11"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.421"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +3.32%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.047"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +13.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.56"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.12"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "142.72"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.416"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.061"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -8.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.544"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +2.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.845.78"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.056"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +8.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08264"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +2.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03024"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.924"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -3.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.18"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +11.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2758"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09226"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7720"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.77"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +5.19%  "
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.63"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +4.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7119"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +3.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.542"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +2.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.132"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08422"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.58"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.267"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.56%  "
